$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 514, shifting rows 514:576 down to 516:578
$ws.Rows("514:515").Insert()

# New row 514: fresh data point dated 45154, quality "Primera", origin "Provincia de Diguillín"
$ws.Cells.Item(514, 1).Value = 7
$ws.Cells.Item(514, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(514, 3).Value = "Ñuble"
$ws.Cells.Item(514, 4).Value = 45154
$ws.Cells.Item(514, 5).Value = 16
$ws.Cells.Item(514, 6).Value = 100112023
$ws.Cells.Item(514, 7).Value = "Brócoli"
$ws.Cells.Item(514, 8).Value = "Sin especificar"
$ws.Cells.Item(514, 9).Value = "Primera"
$ws.Cells.Item(514, 10).Value = 300
$ws.Cells.Item(514, 11).Value = 1000
$ws.Cells.Item(514, 12).Value = 1000
$ws.Cells.Item(514, 13).Value = 1000
$ws.Cells.Item(514, 14).Value = "`$/unidad"
$ws.Cells.Item(514, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(514, 16).Value = 1000
$ws.Cells.Item(514, 17).Value = 1
$ws.Cells.Item(514, 18).Value = "Hortaliza"

# New row 515: fresh data point dated 45154, quality "Segunda", origin "Provincia de Diguillín"
$ws.Cells.Item(515, 1).Value = 7
$ws.Cells.Item(515, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(515, 3).Value = "Ñuble"
$ws.Cells.Item(515, 4).Value = 45154
$ws.Cells.Item(515, 5).Value = 16
$ws.Cells.Item(515, 6).Value = 100112023
$ws.Cells.Item(515, 7).Value = "Brócoli"
$ws.Cells.Item(515, 8).Value = "Sin especificar"
$ws.Cells.Item(515, 9).Value = "Segunda"
$ws.Cells.Item(515, 10).Value = 300
$ws.Cells.Item(515, 11).Value = 800
$ws.Cells.Item(515, 12).Value = 800
$ws.Cells.Item(515, 13).Value = 800
$ws.Cells.Item(515, 14).Value = "`$/unidad"
$ws.Cells.Item(515, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(515, 16).Value = 800
$ws.Cells.Item(515, 17).Value = 1
$ws.Cells.Item(515, 18).Value = "Hortaliza"
